$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 102 (shifts existing rows 102..220 down to 103..221)
$ws.Rows.Item(102).Insert()

# Populate the newly inserted row 102 with data (same as the record that used to sit
# there, but with an updated date and volume figure)
$ws.Cells.Item(102, 1).Value = 4
$ws.Cells.Item(102, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(102, 3).Value = "Los Lagos"
$ws.Cells.Item(102, 4).Value = Get-Date -Year 2022 -Month 3 -Day 18 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(102, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(102, 5).Value = 10
$ws.Cells.Item(102, 6).Value = 100112017
$ws.Cells.Item(102, 7).Value = "Apio"
$ws.Cells.Item(102, 8).Value = "Americana (o)"
$ws.Cells.Item(102, 9).Value = "Primera"
$ws.Cells.Item(102, 10).Value = 50
$ws.Cells.Item(102, 11).Value = 12000
$ws.Cells.Item(102, 12).Value = 12000
$ws.Cells.Item(102, 13).Value = 12000
$ws.Cells.Item(102, 14).Value = "`$/docena de matas"
$ws.Cells.Item(102, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(102, 16).Value = 2000
$ws.Cells.Item(102, 17).Value = 6
$ws.Cells.Item(102, 18).Value = "Hortaliza"
